$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values that look like plain decimal numbers need a quote-prefix so
# Excel keeps storing them as text (matching the source inlineStr cells)
# rather than silently converting them to floating-point numbers; the style
# is then reset back to Normal so no visible formatting change remains.

$ws.Range("D2").Value = '26.260.71'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '1.590.93'
$ws.Range("E3").Value = '  +0.69%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = '''212.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.48%  '
$ws.Range("D6").Value = '''0.502'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.22%  '
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("E8").Value = '  +0.36%  '
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").Value = '''19.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").Value = '''0.0849'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("D12").Value = '1.813.55'
$ws.Range("E12").Value = '  +0.60%  '
$ws.Range("D13").Value = '1.592.05'
$ws.Range("E13").Value = '  +0.84%  '
$ws.Range("E14").Value = '  -0.08%  '
$ws.Range("E15").Value = '  +1.37%  '
$ws.Range("D16").Value = '''64.36'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").Value = '26.257.60'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("D19").Value = '''7.44'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.62%  '
$ws.Range("D20").Value = '''213.49'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.94%  '
$ws.Range("E22").Value = '  +0.66%  '
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("E24").Value = '  -2.50%  '
$ws.Range("D25").Value = '''144.68'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("E26").Value = '  -0.14%  '
$ws.Range("D27").Value = '''7.05'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.71%  '
$ws.Range("E28").Value = '  -0.64%  '
$ws.Range("D29").Value = '''15.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.10%  '
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D31").Value = '''1.15'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.81%  '
$ws.Range("D32").Value = '''3.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").Value = '1.336.82'
$ws.Range("E34").Value = '  +4.70%  '
$ws.Range("E35").Value = '  -1.01%  '
$ws.Range("E36").Value = '  -0.57%  '
$ws.Range("D37").Value = '''0.591'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.08%  '
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("D39").Value = '''0.819'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.36%  '
$ws.Range("E40").Value = '  -4.98%  '
$ws.Range("E42").Value = '  +3.07%  '
$ws.Range("E43").Value = '  -0.07%  '
$ws.Range("D44").Value = '''0.765'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.20%  '
$ws.Range("D45").Value = '1.725.35'
$ws.Range("E45").Value = '  +0.51%  '
$ws.Range("D46").Value = '''61.82'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.97%  '
$ws.Range("D47").Value = '''85.96'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.39%  '
$ws.Range("E48").Value = '  -3.75%  '
$ws.Range("E49").Value = '  -0.60%  '
$ws.Range("D50").Value = '''0.0974'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.89%  '
$ws.Range("D51").Value = '''0.999'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.34%  '
